$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.684.34'
$ws.Range("E2").Value = '  -3.30%  '
$ws.Range("D3").Value = '2.097.63'
$ws.Range("E3").Value = '  -2.20%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").Value = '344.03'
$ws.Range("E5").Value = '  -2.46%  '
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").Value = '0.5138'
$ws.Range("E7").Value = '  -2.78%  '
$ws.Range("D8").Value = '0.4408'
$ws.Range("E8").Value = '  -3.46%  '
$ws.Range("D9").Value = '52.76'
$ws.Range("E9").Value = '  -2.99%  '
$ws.Range("D10").Value = '0.09158'
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").Value = '24.98'
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").Value = '2.108.97'
$ws.Range("E13").Value = '  -1.42%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '8.257'
$ws.Range("E14").Value = '  +1.25%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '6.760'
$ws.Range("E15").Value = '  -2.23%  '
$ws.Range("D16").Value = '99.64'
$ws.Range("E16").Value = '  -2.63%  '
$ws.Range("E17").Value = '  -2.48%  '
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Value = '20.79'
$ws.Range("E19").Value = '  +6.23%  '
$ws.Range("D20").Value = '0.06626'
$ws.Range("E20").Value = '  -1.53%  '
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Value = '6.187'
$ws.Range("E22").Value = '  -3.29%  '
$ws.Range("D23").Value = '29.751.08'
$ws.Range("E23").Value = '  -3.39%  '
$ws.Range("D24").Value = '12.59'
$ws.Range("E24").Value = '  -2.53%  '
$ws.Range("D25").Value = '2.316'
$ws.Range("E25").Value = '  -3.06%  '
$ws.Range("D26").Value = '2.352.10'
$ws.Range("E26").Value = '  -1.53%  '
$ws.Range("E27").Value = '  -3.19%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '2.529'
$ws.Range("E28").Value = '  -2.20%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '162.04'
$ws.Range("E29").Value = '  -1.84%  '
$ws.Range("D30").Value = '132.82'
$ws.Range("E30").Value = '  -3.28%  '
$ws.Range("D31").Value = '1.130'
$ws.Range("E31").Value = '  -6.94%  '
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").Value = '1.662'
$ws.Range("E32").Value = '  -1.65%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = '0.1049'
$ws.Range("E33").Value = '  -3.56%  '
$ws.Range("D34").Value = '6.167'
$ws.Range("E34").Value = '  -3.98%  '
$ws.Range("D35").Value = '3.942'
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '10.44'
$ws.Range("E36").Value = '  +0.95%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '6.025'
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("D38").Value = '0.02564'
$ws.Range("D39").Value = '0.06723'
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("D40").Value = '12.44'
$ws.Range("E40").Value = '  -1.62%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6861'
$ws.Range("E41").Value = '  -1.45%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.2232'
$ws.Range("E42").Value = '  -4.40%  '
$ws.Range("D43").Value = '1.290'
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("D44").Value = '0.6661'
$ws.Range("E44").Value = '  +2.77%  '
$ws.Range("D45").Value = '14.21'
$ws.Range("E45").Value = '  -4.53%  '
$ws.Range("E46").Value = '  -2.38%  '
$ws.Range("D47").Value = '3.609'
$ws.Range("E47").Value = '  -4.08%  '
$ws.Range("D48").Value = '0.00000000349'
$ws.Range("E48").Value = '  -5.34%  '
$ws.Range("D49").Value = '1.220'
$ws.Range("D50").Value = '82.35'
$ws.Range("E50").Value = '  -1.43%  '
$ws.Range("D51").Value = '0.3324'
$ws.Range("E51").Value = '  -2.80%  '
